$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = 6233
$ws.Range("C21").Value = 989
$ws.Range("D21").Value = 5610346
$ws.Range("E21").Value = 900.1036419059843
$ws.Range("F21").Value = 8.193022044783893
$ws.Range("G21").Value = 4.324894514767941
$ws.Range("H21").Value = 28.0616721585986
